$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "all"
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "saturdays"

[void]$ws.Range("A62").Select()
[void]$ws2.Range("D27").Select()
[void]$ws2.Activate()
Write-Host "done"
